# Update the "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, reflecting a refreshed data scrape (gh-pages output
# regenerated at commit 456a3b4).
#
# Sheet order in the workbook:
#   1 -> 展览     (exhibitions)
#   2 -> 演出     (performances)
#   3 -> 本地生活 (local life)
#   4 -> 全部类型 (all types / combined)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- 展览 (sheet 1) ---
$ws1.Range("F2").Value = 1899
$ws1.Range("F3").Value = 1525
$ws1.Range("F4").Value = 890
$ws1.Range("F5").Value = 788
$ws1.Range("F6").Value = 13392
$ws1.Range("F7").Value = 13258
$ws1.Range("F9").Value = 779
$ws1.Range("F11").Value = 568
$ws1.Range("F13").Value = 689
$ws1.Range("F23").Value = 764
$ws1.Range("F24").Value = 21

# --- 演出 (sheet 2) ---
$ws2.Range("F2").Value = 97
$ws2.Range("F7").Value = 129
$ws2.Range("F8").Value = 12

# --- 本地生活 (sheet 3) ---
$ws3.Range("F2").Value = 197

# --- 全部类型 (sheet 4) ---
$ws4.Range("F2").Value = 197
$ws4.Range("F3").Value = 1899
$ws4.Range("F4").Value = 1525
$ws4.Range("F5").Value = 890
$ws4.Range("F6").Value = 97
$ws4.Range("F7").Value = 788
$ws4.Range("F8").Value = 13392
$ws4.Range("F9").Value = 13258
$ws4.Range("F11").Value = 779
$ws4.Range("F13").Value = 568
$ws4.Range("F15").Value = 689
$ws4.Range("F30").Value = 764
$ws4.Range("F31").Value = 129
$ws4.Range("F32").Value = 12
$ws4.Range("F33").Value = 21

$wb.Save()
